$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-25 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-26 Friday", 2)

$d.Content.Find.Execute("28×59=1652", $true, $false, $false, $false, $false, $true, 1, $false, "60×49=2940", 2)
$d.Content.Find.Execute("84×17=1428", $true, $false, $false, $false, $false, $true, 1, $false, "24×64=1536", 2)
$d.Content.Find.Execute("78×47=3666", $true, $false, $false, $false, $false, $true, 1, $false, "15×63=945", 2)
$d.Content.Find.Execute("33×65=2145", $true, $false, $false, $false, $false, $true, 1, $false, "76×34=2584", 2)
$d.Content.Find.Execute("41×76=3116", $true, $false, $false, $false, $false, $true, 1, $false, "83×26=2158", 2)
$d.Content.Find.Execute("21×97=2037", $true, $false, $false, $false, $false, $true, 1, $false, "83×67=5561", 2)
$d.Content.Find.Execute("61×98=5978", $true, $false, $false, $false, $false, $true, 1, $false, "51×68=3468", 2)
$d.Content.Find.Execute("65×30=1950", $true, $false, $false, $false, $false, $true, 1, $false, "72×77=5544", 2)
$d.Content.Find.Execute("91×91=8281", $true, $false, $false, $false, $false, $true, 1, $false, "78×92=7176", 2)
$d.Content.Find.Execute("93×32=2976", $true, $false, $false, $false, $false, $true, 1, $false, "39×29=1131", 2)
$d.Content.Find.Execute("40×21=840", $true, $false, $false, $false, $false, $true, 1, $false, "82×25=2050", 2)
$d.Content.Find.Execute("95×49=4655", $true, $false, $false, $false, $false, $true, 1, $false, "61×49=2989", 2)
$d.Content.Find.Execute("86×31=2666", $true, $false, $false, $false, $false, $true, 1, $false, "22×29=638", 2)
$d.Content.Find.Execute("11×65=715", $true, $false, $false, $false, $false, $true, 1, $false, "30×61=1830", 2)
$d.Content.Find.Execute("66×37=2442", $true, $false, $false, $false, $false, $true, 1, $false, "94×12=1128", 2)
$d.Content.Find.Execute("84×49=4116", $true, $false, $false, $false, $false, $true, 1, $false, "75×80=6000", 2)
$d.Content.Find.Execute("76×78=5928", $true, $false, $false, $false, $false, $true, 1, $false, "82×24=1968", 2)
$d.Content.Find.Execute("76×96=7296", $true, $false, $false, $false, $false, $true, 1, $false, "64×56=3584", 2)
$d.Content.Find.Execute("26×60=1560", $true, $false, $false, $false, $false, $true, 1, $false, "24×71=1704", 2)
$d.Content.Find.Execute("55×15=825", $true, $false, $false, $false, $false, $true, 1, $false, "73×82=5986", 2)
$d.Content.Find.Execute("94×51=4794", $true, $false, $false, $false, $false, $true, 1, $false, "42×69=2898", 2)
$d.Content.Find.Execute("72×23=1656", $true, $false, $false, $false, $false, $true, 1, $false, "52×24=1248", 2)
$d.Content.Find.Execute("25×15=375", $true, $false, $false, $false, $false, $true, 1, $false, "32×33=1056", 2)
$d.Content.Find.Execute("60×47=2820", $true, $false, $false, $false, $false, $true, 1, $false, "92×89=8188", 2)
$d.Content.Find.Execute("15×64=960", $true, $false, $false, $false, $false, $true, 1, $false, "58×15=870", 2)
